# Update the "Förändrad" (Changed) date column (C) for all data rows.
# The workbook was re-generated on a later day, so every row's "last
# changed/checked" timestamp in column C moves from serial date 45203
# (2023-10-04) to serial date 45204 (2023-10-05).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp
if ($lastRow -lt 2) { $lastRow = 2 }

$range = $ws.Range("C2:C$lastRow")
$range.Value = 45204
